# Commit: "change Ranker to ranker 2"
# Rename the "Ranker.com" domain (and the stray "headerbiddingproptest.com"
# domain on row 21, which was a leftover copy/paste bug) to lowercase
# "ranker.com" across the Ranker test-case block (rows 13-23):
#   - Column D (DOMAIN display column) -> "ranker.com"
#   - Column E (BODY JSON) -> the "domain" field inside the JSON payload
#     is renamed to "ranker.com" as well (row 17's JSON has no "domain"
#     field at all, so it is left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = 13..23

foreach ($r in $rows) {
    # Column D: DOMAIN column always becomes "ranker.com"
    $ws.Cells.Item($r, 4).Value2 = "ranker.com"

    # Column E: BODY column - rewrite the "domain" value inside the JSON text
    $body = $ws.Cells.Item($r, 5).Value2
    if ($body -ne $null) {
        $newBody = $body.Replace('"domain": "Ranker.com"', '"domain": "ranker.com"')
        $newBody = $newBody.Replace('"domain": "headerbiddingproptest.com"', '"domain": "ranker.com"')
        $ws.Cells.Item($r, 5).Value2 = $newBody
    }
}

# Row 21's JSON body got shorter (now mirrors the shape of row 19's body),
# so Excel's autofit shrank the row height from 270.75 to 242.25.
$ws.Rows.Item(21).RowHeight = 242.25

# The saved file also carries an updated view/selection: Excel had scrolled
# back up and the user had clicked on E23 instead of A23.
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Range("E23").Select()
